# Fruta / hortaliza, semanal
# Insert two new daily-price rows for "Choclo" (Macroferia Regional de Talca)
# right after the existing row 87, pushing the remainder of the table down by
# two rows (old row 88 -> 90, ..., old row 160 -> 162) and growing the used
# range from A1:R160 to A1:R162.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new records.
$ws.Rows("88:89").Insert()

# New row 88: "Primera" quality, 50000 units @ 250.
$row88 = New-Object 'object[,]' 1,17
$row88[0,0]  = "Macroferia Regional de Talca"  # B88
$row88[0,1]  = "Maule"                          # C88
$row88[0,2]  = 44566                            # D88 (Fecha)
$row88[0,3]  = 7                                # E88 (Codreg)
$row88[0,4]  = 100112024                        # F88 (Categoría ID)
$row88[0,5]  = "Choclo"                         # G88
$row88[0,6]  = "Choclero"                       # H88 (Variedad)
$row88[0,7]  = "Primera"                        # I88 (Calidad)
$row88[0,8]  = 50000                            # J88 (Volumen)
$row88[0,9]  = 250                              # K88 (Precio mínimo)
$row88[0,10] = 250                              # L88 (Precio máximo)
$row88[0,11] = 250                              # M88 (Precio promedio ponderado)
$row88[0,12] = "`$/unidad"                      # N88 (Unidad de comercialización)
$row88[0,13] = "Región del Maule"               # O88 (Origen)
$row88[0,14] = 250                              # P88 (Precio $/Kg)
$row88[0,15] = 1                                # Q88 (Kg o Unidades)
$row88[0,16] = "Hortaliza"                      # R88 (Clasificación)

$ws.Range("A88").Value = 5
$ws.Range("B88:R88").Value = $row88

# New row 89: "Segunda" quality, 8000 units @ 100.
$row89 = New-Object 'object[,]' 1,17
$row89[0,0]  = "Macroferia Regional de Talca"  # B89
$row89[0,1]  = "Maule"                          # C89
$row89[0,2]  = 44566                            # D89 (Fecha)
$row89[0,3]  = 7                                # E89 (Codreg)
$row89[0,4]  = 100112024                        # F89 (Categoría ID)
$row89[0,5]  = "Choclo"                         # G89
$row89[0,6]  = "Choclero"                       # H89 (Variedad)
$row89[0,7]  = "Segunda"                        # I89 (Calidad)
$row89[0,8]  = 8000                             # J89 (Volumen)
$row89[0,9]  = 100                              # K89 (Precio mínimo)
$row89[0,10] = 100                              # L89 (Precio máximo)
$row89[0,11] = 100                              # M89 (Precio promedio ponderado)
$row89[0,12] = "`$/unidad"                      # N89 (Unidad de comercialización)
$row89[0,13] = "Región del Maule"               # O89 (Origen)
$row89[0,14] = 100                              # P89 (Precio $/Kg)
$row89[0,15] = 1                                # Q89 (Kg o Unidades)
$row89[0,16] = "Hortaliza"                      # R89 (Clasificación)

$ws.Range("A89").Value = 5
$ws.Range("B89:R89").Value = $row89
